$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new top data row for 2022-Q4
#    and shift the existing quarters down, renumbering the index column.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Range("A2:D2").EntireRow.Insert()

# new row 2 had formatting stripped on B:D (same as other data rows); give A2
# back the index-column style used by the sibling rows, then set the values.
$summary.Range("B2:D2").ClearFormats()
$summary.Range("A3").Copy($summary.Range("A2"))

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.01

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# ---------------------------------------------------------------------------
# 2) Insert the new "2022-Q4" fund-holdings sheet right after "总计", copying
#    the layout/style of the "2021-Q3" sheet and replacing its data with the
#    2022-Q4 figures.
# ---------------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2021-Q3")
$templateSheet.Copy($null, $summary)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# drop the template's extra data rows, keep header (row1) + 2 data rows
$q4.Range("A4:H11").EntireRow.Delete()

$q4.Range("D1").Value = "基金规模"

$q4.Range("B2:G3").NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "006658"
$q4.Range("C2").Value = "财通中证香港红利等权投资指数A"
$q4.Range("D2").Value = "0.14"
$q4.Range("E2").Value = "89.84"
$q4.Range("F2").Value = "3.06"
$q4.Range("G2").Value = "0.0043"
$q4.Range("H2").Value = 10

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "006659"
$q4.Range("C3").Value = "财通中证香港红利等权投资指数C"
$q4.Range("D3").Value = "0.04"
$q4.Range("E3").Value = "89.84"
$q4.Range("F3").Value = "3.06"
$q4.Range("G3").Value = "0.0012"
$q4.Range("H3").Value = 10
